$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.662.00"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.775.31"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "664.70"
$ws.Range("E5").Value = "  +6.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.01"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").Value = "3.773.89"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.99"
$ws.Range("E12").Value = "  +5.85%  "
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.23"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "4.413.22"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "3.774.65"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "69.485.76"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.68"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.10"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.40"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.73"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.20"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "3.920.85"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.27"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.91"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +18.47%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "3.727.88"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.98"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.958"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "45.52"
$ws.Range("E45").Value = "  +7.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.90"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.97"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.299"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.40"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  +0.72%  "
